# MCU_Selection.xlsx edit script
# - Update the "Parts" sheet (rename to "Microcontroller"), fix a couple of data cells
# - Insert a brand-new "ADC" sheet between it and the README sheet
# - Rename "README" to "Readme"
# - Fix up sheet views / active tab to match the authored state

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename the Parts sheet and insert the new ADC sheet right after it.
#    NOTE: Worksheets.Add() shifts sheet indices around, which can leave
#    previously-captured sheet variables pointing at the wrong sheet, so
#    every sheet reference we keep past the Add() call below is
#    (re-)fetched by its stable name afterwards.
# ------------------------------------------------------------------
$wsMcu = $wb.Worksheets.Item(1)
$wsMcu.Name = "Microcontroller"

$wsAdc = $wb.Worksheets.Add($null, $wsMcu)
$wsAdc.Name = "ADC"

# Re-fetch by name now that the sheet collection has shifted around.
$wsMcu = $wb.Worksheets.Item("Microcontroller")
$wsAdc = $wb.Worksheets.Item("ADC")
$wsReadme = $wb.Worksheets.Item("README")
$wsReadme.Name = "Readme"
$wsReadme = $wb.Worksheets.Item("Readme")

# ------------------------------------------------------------------
# 2. Fix a couple of data points on the Microcontroller sheet
# ------------------------------------------------------------------
# Row 3 was a duplicate "ATMega328p" entry - it should actually be the
# ATMega168p part.
$wsMcu.Range("A3").Value = "ATMega168p"

# Digikey stock count for the ATSAMD21J15B row dropped to 0.
$wsMcu.Range("K9").Value = 0

# ------------------------------------------------------------------
# 3. Populate the new "ADC" sheet
# ------------------------------------------------------------------

# Header row
$wsAdc.Range("A1").Value = "ADC"
$wsAdc.Range("B1").Value = "Manufacturer"
$wsAdc.Range("C1").Value = "Res Bits"
$wsAdc.Range("D1").Value = "# Chan"
$wsAdc.Range("E1").Value = "Sample Rate"
$wsAdc.Range("F1").Value = "Interface"
$wsAdc.Range("G1").Value = "Price ($)"
$wsAdc.Range("G1").NumberFormat = """$""#,##0.00"
$wsAdc.Range("H1").Value = "Digikey Stock"
$wsAdc.Range("I1").Value = "Link"
$wsAdc.Range("J1").Value = "Additional notes"

# Row 2 - ADC108S102
$wsAdc.Range("A2").Value = "ADC108S102"
$wsAdc.Range("B2").Value = "Texas Inst"
$wsAdc.Range("C2").Value = 10
$wsAdc.Range("D2").Value = 8
$wsAdc.Range("E2").Value = "500k-1M"
$wsAdc.Range("F2").Value = "SPI"
$wsAdc.Range("G2").Value = 8.33
$wsAdc.Range("G2").NumberFormat = """$""#,##0.00"
$wsAdc.Range("H2").Value = 173
$wsAdc.Range("I2").Value = "https://www.digikey.com/en/products/detail/texas-instruments/ADC108S102CIMT-NOPB/953338 "
$wsAdc.Hyperlinks.Add($wsAdc.Range("I2"), "https://www.digikey.com/en/products/detail/texas-instruments/ADC108S102CIMT-NOPB/953338") | Out-Null
$wsAdc.Range("I2").Style = "Hyperlink"

# Row 3 - ADC128S102
$wsAdc.Range("A3").Value = "ADC128S102"
$wsAdc.Range("B3").Value = "Texas Inst"
$wsAdc.Range("C3").Value = 12
$wsAdc.Range("D3").Value = 8
$wsAdc.Range("E3").Value = "500k-1M"
$wsAdc.Range("F3").Value = "SPI"
$wsAdc.Range("G3").Value = 8.33
$wsAdc.Range("G3").NumberFormat = """$""#,##0.00"
$wsAdc.Range("H3").Value = 3708
$wsAdc.Range("I3").Value = "https://www.digikey.com/en/products/detail/texas-instruments/ADC128S102CIMTX-NOPB/1870710 "
$wsAdc.Hyperlinks.Add($wsAdc.Range("I3"), "https://www.digikey.com/en/products/detail/texas-instruments/ADC128S102CIMTX-NOPB/1870710") | Out-Null
$wsAdc.Range("I3").Style = "Hyperlink"

# Row 4 - ADS7955SDBT
$wsAdc.Range("A4").Value = "ADS7955SDBT"
$wsAdc.Range("B4").Value = "Texas Inst"
$wsAdc.Range("C4").Value = 10
$wsAdc.Range("D4").Value = 8
$wsAdc.Range("E4").Value = "1M"
$wsAdc.Range("F4").Value = "SPI"
$wsAdc.Range("G4").Value = 6.66
$wsAdc.Range("G4").NumberFormat = """$""#,##0.00"
$wsAdc.Range("H4").Value = 422
$wsAdc.Range("I4").Value = "https://www.digikey.com/en/products/detail/texas-instruments/ADS7955SDBT/1880865 "
$wsAdc.Hyperlinks.Add($wsAdc.Range("I4"), "https://www.digikey.com/en/products/detail/texas-instruments/ADS7955SDBT/1880865") | Out-Null
$wsAdc.Range("I4").Style = "Hyperlink"

# Column widths to roughly match the authored layout
$wsAdc.Columns.Item(1).ColumnWidth = 13.7
$wsAdc.Columns.Item(2).ColumnWidth = 11.17
$wsAdc.Columns.Item(5).ColumnWidth = 10.87
$wsAdc.Columns.Item(8).ColumnWidth = 14.17
$wsAdc.Columns.Item(9).ColumnWidth = 13.52

$wsAdc.PageSetup.Orientation = 1

# ------------------------------------------------------------------
# 4. Sheet view / selection bookkeeping, last Activate() wins the
#    workbook's active tab.
# ------------------------------------------------------------------
$wsMcu.Activate()
$wsMcu.Range("G6").Select()

$wsReadme.Activate()
$wsReadme.Range("F14").Select()

$wsAdc.Activate()
$wsAdc.Range("C13").Select()
